$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: new "logs" (0xD) command definition row ---
# A16 / D16 take the plain white-fill, no-alignment style already used by
# cells such as G14/H14 (fillId4, no special alignment).
$ws.Range("G14").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("D16").PasteSpecial(-4122)

# C16 and E16:I16 take the plain white-fill centered style used by
# D15/H15/I15; E16:I16 later get merged into a single cell holding the
# JSON description.
$ws.Range("D15").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("E16:I16").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$ws.Range("A16").Value = "logs"
$ws.Range("C16").Value = ">0x002"
$ws.Range("D16").Value = "-"
$ws.Range("E16").Value = 'JSON { logs:[ {"severity": number, "msg": string} ] }'

$ws.Range("E16:I16").Merge()
